# Replace every occurrence of the literal string "(0, 0)" in column F
# (the "EXP Start Point" column) with "(nan, nan)" across all worksheets
# in the workbook.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

    for ($r = 1; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 6)  # column F
        if ($cell.Value2 -eq "(0, 0)") {
            $cell.Value2 = "(nan, nan)"
        }
    }
}
